$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The five "TabQuery"/"StatQuery" SQL cells (C2, B2, B3, B4, B5, B6, B7) all share the
# same buggy JOIN block:
#   df_participant prt ON std.id = prt."study.id"
#   df_diagnoses   dgn ON prt.id = dgn."participant.id"
#   df_treatments  trt ON prt.id = trt."participant.id"
#   df_treatment_resp trr ON prt.id = trr."participant.id"
#   df_survival    srv ON prt.id = srv."participant.id"
#   df_reference_files rfs ON std.id = rfs."study.id"
# Update every occurrence to join on the real key columns (study_id / participant_id).
$oldPattern = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newPattern = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $queryCells) {
    $cell = $ws.Range($addr)
    $val = $cell.Value()
    if ($val -ne $null -and $val.Contains($oldPattern)) {
        $cell.Value = $val.Replace($oldPattern, $newPattern)
    }
}

# Column C was auto-fit ("bestFit") at ~60.83 chars; widen it to a fixed 71 chars
# (matches width="71" customWidth="1" with bestFit removed in the target file).
$ws.Columns.Item(3).ColumnWidth = 70.2

# Scroll the sheet so row 6 is at the top of the view (matches topLeftCell="A6"),
# keeping the original active-cell selection at C7.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
